$wb = $excel.ActiveWorkbook

# Rename "temperature_c" sheet to "temperature"
$wsTemp = $wb.Worksheets.Item("temperature_c")
$wsTemp.Name = "temperature"

# Select the "temperature" sheet and set its active cell/selection
$wsTemp.Select()
$wsTemp.Range("Q34").Select()

# Select the "genotype" sheet and set its active cell/selection
$wsGeno = $wb.Worksheets.Item("genotype")
$wsGeno.Select()
$wsGeno.Range("J19:J20").Select()

# Make "temperature" the final active sheet (so activeTab points at it)
$wsTemp.Select()
